# Add payment 71277628 (Cash) 2025-08-18T16:53:13
#
# Prior to this edit, row 20's "phone" cell (A20) had been written as a
# literal/text "79174445" instead of the numeric value used by every other
# row. This edit both normalizes A20 back to a number and appends the new
# payment as row 21 (mirroring the existing column layout: phone, amount,
# method, timestamp, original_amount, discount_applied, final_amount,
# birthday_discount, points_redeemed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 20: A20 should be the number 79174445, not text ---
$ws.Cells.Item(20, 1).Value = 79174445

# --- Append row 21: new Cash payment for phone 71277628 ---
# Column A (phone) keeps its leading-apostrophe text form, consistent with
# how this sheet records a freshly-ingested phone number before it gets
# normalized to a number (see A20's prior state, above).
$ws.Cells.Item(21, 1).Value = "'71277628"
$ws.Cells.Item(21, 3).Value = "Cash"
$ws.Cells.Item(21, 4).Value = "2025-08-18T16:53:13"
$ws.Cells.Item(21, 5).Value = 80
$ws.Cells.Item(21, 7).Value = 80
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
